# The sheet's 4 data columns (A=deaths, B=date, C=hospitalized, D=total)
# are being reordered into (A=total, B=hospitalized, C=deaths, D=date),
# i.e. column A swaps with column D, and column B swaps with column C.
# We do this with a pair of column-level Cut/Paste moves via a scratch
# area (columns F:I) so that per-cell number formatting (the date style
# used in the old "date" column) travels together with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 23

# Stage 1: move the whole A1:D23 block out to a scratch area (F1:I23)
# so F=old A, G=old B, H=old C, I=old D.
$ws.Range("A1:D23").Cut($ws.Range("F1:I23")) | Out-Null

# Make sure the original area is completely empty (including formatting)
# before we rebuild it.
$ws.Range("A1:D23").Clear() | Out-Null

# Stage 2: move each scratch column into its new home column.
# New A = old D (I)    New B = old C (H)
# New C = old A (F)    New D = old B (G)
$ws.Range("I1:I" + $lastRow).Cut($ws.Range("A1:A" + $lastRow)) | Out-Null
$ws.Range("H1:H" + $lastRow).Cut($ws.Range("B1:B" + $lastRow)) | Out-Null
$ws.Range("F1:F" + $lastRow).Cut($ws.Range("C1:C" + $lastRow)) | Out-Null
$ws.Range("G1:G" + $lastRow).Cut($ws.Range("D1:D" + $lastRow)) | Out-Null

# Remove any leftovers in the scratch area.
$ws.Range("F1:I23").Clear() | Out-Null

# The Cut/Paste moves above can leave behind empty-but-present cells
# (no value, but still serialized) wherever a column used to have data
# in a row that is now blank for that column. Strip those out so blank
# cells are fully absent, matching a freshly blank cell.
for ($r = 1; $r -le $lastRow; $r++) {
    for ($colIx = 1; $colIx -le 4; $colIx++) {
        $cell = $ws.Cells.Item($r, $colIx)
        if ($cell.Value2 -eq $null) {
            $cell.Clear() | Out-Null
        }
    }
}

# Select the whole table, matching the new sheet selection state.
$ws.Range("A1:D23").Select() | Out-Null
